$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.228.54"
$ws.Range("D2").Style = $origStyle
$origStyle = $ws.Range("E2").Style
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("E2").Style = $origStyle

$origStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.631.04"
$ws.Range("D3").Style = $origStyle
$origStyle = $ws.Range("E3").Style
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("E3").Style = $origStyle

$origStyle = $ws.Range("E4").Style
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E4").Style = $origStyle

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.08"
$ws.Range("D5").Style = $origStyle
$origStyle = $ws.Range("E5").Style
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("E5").Style = $origStyle

$origStyle = $ws.Range("E6").Style
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.63%  "
$ws.Range("E6").Style = $origStyle

$origStyle = $ws.Range("E7").Style
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E7").Style = $origStyle

$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.256"
$ws.Range("D8").Style = $origStyle
$origStyle = $ws.Range("E8").Style
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("E8").Style = $origStyle

$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0625"
$ws.Range("D9").Style = $origStyle
$origStyle = $ws.Range("E9").Style
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("E9").Style = $origStyle

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.21"
$ws.Range("D10").Style = $origStyle
$origStyle = $ws.Range("E10").Style
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("E10").Style = $origStyle

$origStyle = $ws.Range("E11").Style
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("E11").Style = $origStyle

$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.637.78"
$ws.Range("D12").Style = $origStyle
$origStyle = $ws.Range("E12").Style
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.96%  "
$ws.Range("E12").Style = $origStyle

$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.12"
$ws.Range("D13").Style = $origStyle
$origStyle = $ws.Range("E13").Style
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("E13").Style = $origStyle

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.543"
$ws.Range("D14").Style = $origStyle
$origStyle = $ws.Range("E14").Style
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("E14").Style = $origStyle

$origStyle = $ws.Range("B15").Style
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "Litecoin"
$ws.Range("B15").Style = $origStyle
$origStyle = $ws.Range("C15").Style
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("C15").Style = $origStyle
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.91"
$ws.Range("D15").Style = $origStyle
$origStyle = $ws.Range("E15").Style
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.47%  "
$ws.Range("E15").Style = $origStyle

$origStyle = $ws.Range("B16").Style
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("B16").Style = $origStyle
$origStyle = $ws.Range("C16").Style
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("C16").Style = $origStyle
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.185.06"
$ws.Range("D16").Style = $origStyle
$origStyle = $ws.Range("E16").Style
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("E16").Style = $origStyle

$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₃0734"
$ws.Range("D17").Style = $origStyle
$origStyle = $ws.Range("E17").Style
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.02%  "
$ws.Range("E17").Style = $origStyle

$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "216.77"
$ws.Range("D18").Style = $origStyle
$origStyle = $ws.Range("E18").Style
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.75%  "
$ws.Range("E18").Style = $origStyle

$origStyle = $ws.Range("E19").Style
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("E19").Style = $origStyle

$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.94"
$ws.Range("D20").Style = $origStyle
$origStyle = $ws.Range("E20").Style
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.09%  "
$ws.Range("E20").Style = $origStyle

$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.39"
$ws.Range("D21").Style = $origStyle
$origStyle = $ws.Range("E21").Style
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.46%  "
$ws.Range("E21").Style = $origStyle

$origStyle = $ws.Range("E22").Style
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.72%  "
$ws.Range("E22").Style = $origStyle

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.12"
$ws.Range("D23").Style = $origStyle
$origStyle = $ws.Range("E23").Style
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.28%  "
$ws.Range("E23").Style = $origStyle

$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.22"
$ws.Range("D24").Style = $origStyle
$origStyle = $ws.Range("E24").Style
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("E24").Style = $origStyle

$origStyle = $ws.Range("E25").Style
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E25").Style = $origStyle

$origStyle = $ws.Range("E26").Style
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.72%  "
$ws.Range("E26").Style = $origStyle

$origStyle = $ws.Range("E27").Style
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("E27").Style = $origStyle

$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.58"
$ws.Range("D28").Style = $origStyle
$origStyle = $ws.Range("E28").Style
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.83%  "
$ws.Range("E28").Style = $origStyle

$origStyle = $ws.Range("E29").Style
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("E29").Style = $origStyle

$origStyle = $ws.Range("E30").Style
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.06%  "
$ws.Range("E30").Style = $origStyle

$origStyle = $ws.Range("E31").Style
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("E31").Style = $origStyle

$origStyle = $ws.Range("E32").Style
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.99%  "
$ws.Range("E32").Style = $origStyle

$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.322.13"
$ws.Range("D33").Style = $origStyle
$origStyle = $ws.Range("E33").Style
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.84%  "
$ws.Range("E33").Style = $origStyle

$origStyle = $ws.Range("E34").Style
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.25%  "
$ws.Range("E34").Style = $origStyle

$origStyle = $ws.Range("E35").Style
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("E35").Style = $origStyle

$origStyle = $ws.Range("E36").Style
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.84%  "
$ws.Range("E36").Style = $origStyle

$origStyle = $ws.Range("E37").Style
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.97%  "
$ws.Range("E37").Style = $origStyle

$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.847"
$ws.Range("D38").Style = $origStyle
$origStyle = $ws.Range("E38").Style
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("E38").Style = $origStyle

$origStyle = $ws.Range("E39").Style
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("E39").Style = $origStyle

$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.25"
$ws.Range("D40").Style = $origStyle
$origStyle = $ws.Range("E40").Style
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.09%  "
$ws.Range("E40").Style = $origStyle

$origStyle = $ws.Range("E41").Style
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.25%  "
$ws.Range("E41").Style = $origStyle

$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "63.65"
$ws.Range("D42").Style = $origStyle
$origStyle = $ws.Range("E42").Style
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.13%  "
$ws.Range("E42").Style = $origStyle

$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.767.49"
$ws.Range("D43").Style = $origStyle
$origStyle = $ws.Range("E43").Style
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.66%  "
$ws.Range("E43").Style = $origStyle

$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.74"
$ws.Range("D45").Style = $origStyle

$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.61"
$ws.Range("D46").Style = $origStyle
$origStyle = $ws.Range("E46").Style
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.39%  "
$ws.Range("E46").Style = $origStyle

$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₆0107"
$ws.Range("D47").Style = $origStyle
$origStyle = $ws.Range("E47").Style
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("E47").Style = $origStyle

$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.810"
$ws.Range("D48").Style = $origStyle
$origStyle = $ws.Range("E48").Style
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +21.10%  "
$ws.Range("E48").Style = $origStyle

$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0515"
$ws.Range("D49").Style = $origStyle
$origStyle = $ws.Range("E49").Style
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("E49").Style = $origStyle

$origStyle = $ws.Range("E50").Style
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.33%  "
$ws.Range("E50").Style = $origStyle

$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0957"
$ws.Range("D51").Style = $origStyle
$origStyle = $ws.Range("E51").Style
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.88%  "
$ws.Range("E51").Style = $origStyle
